$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto market data.
# Numeric-looking Price values must be forced to remain text (matching the original
# inline-string cell type) instead of being auto-converted to numbers by Excel.

$ws.Range("D2").Value = '69.545.91'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '3.496.23'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '189.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.212'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.647'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.00'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.19%  '
$ws.Range("E12").Value = '  -4.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").Value = '4.056.59'
$ws.Range("E14").Value = '  -1.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '598.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").Value = '69.629.36'
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.69%  '
$ws.Range("D19").Value = '3.500.18'
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.985'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '105.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("E26").Value = '  +2.13%  '
$ws.Range("E27").Value = '  -2.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.32'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.19%  '
$ws.Range("E30").Value = '  -3.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.65%  '
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("E33").Value = '  -1.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.16'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.06%  '
$ws.Range("D38").Value = '3.628.22'
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("E39").Value = '  -4.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '500.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.51%  '
$ws.Range("E42").Value = '  -4.41%  '
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("E47").Value = '  -4.50%  '
$ws.Range("E48").Value = '  -6.16%  '
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("E51").Value = '  -8.18%  '
